$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant xlPasteFormats (used to copy the cell style from the
# neighboring column K cell into the new column L cell before writing values).
$xlPasteFormats = -4122

function Set-CellWithStyleOfNeighbor {
    param($SourceAddress, $TargetAddress, $Value)
    $ws.Range($SourceAddress).Copy() | Out-Null
    $ws.Range($TargetAddress).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($TargetAddress).Value = $Value
}

# New column L ("2020") values, copying number formatting/style from column K
Set-CellWithStyleOfNeighbor "K4"  "L4"  2020
Set-CellWithStyleOfNeighbor "K5"  "L5"  1.2
Set-CellWithStyleOfNeighbor "K6"  "L6"  1.7
Set-CellWithStyleOfNeighbor "K7"  "L7"  0.4
Set-CellWithStyleOfNeighbor "K8"  "L8"  3.3
Set-CellWithStyleOfNeighbor "K9"  "L9"  3.9
Set-CellWithStyleOfNeighbor "K10" "L10" 2.4
Set-CellWithStyleOfNeighbor "K11" "L11" 95.5
Set-CellWithStyleOfNeighbor "K12" "L12" 94.4
Set-CellWithStyleOfNeighbor "K13" "L13" 97.2

# Scroll the view so column C is the left-most visible column, then select
# the newly populated column L range (matches the saved view state in the diff).
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L4:L13").Select() | Out-Null
